$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" header column (H1), matching the formatting of the existing
# header cells (bold, bordered, centered) by copying G1's format onto H1.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("H1").Value = "Save"

# New data value for the Save column (H2)
$ws.Range("H2").Value = 1
